$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 7).Value = 1.97
$ws.Cells.Item(2, 10).Value = 3.8
$ws.Cells.Item(2, 16).Value = 1.9
$ws.Cells.Item(2, 23).Value = 2.02
$ws.Cells.Item(2, 25).Value = 19
$ws.Cells.Item(2, 28).Value = 9.4
$ws.Cells.Item(2, 29).Value = 970
$ws.Cells.Item(3, 7).Value = 2.76
$ws.Cells.Item(3, 8).Value = 3.15
$ws.Cells.Item(3, 9).Value = 3.9
$ws.Cells.Item(3, 12).Value = 1.41
$ws.Cells.Item(3, 22).Value = 1.35
$ws.Cells.Item(4, 28).Value = 970
$ws.Cells.Item(4, 36).Value = 980
$ws.Cells.Item(6, 12).Value = 1.58
$ws.Cells.Item(7, 6).Value = 2.04
$ws.Cells.Item(7, 8).Value = 3.2
$ws.Cells.Item(7, 11).Value = 4.4
$ws.Cells.Item(7, 12).Value = 1.29
$ws.Cells.Item(7, 14).Value = 4.9
$ws.Cells.Item(7, 15).Value = 1.2
$ws.Cells.Item(7, 17).Value = 1.59
$ws.Cells.Item(7, 19).Value = 2.46
$ws.Cells.Item(7, 20).Value = 1.55
$ws.Cells.Item(7, 24).Value = 28
$ws.Cells.Item(7, 25).Value = 22
$ws.Cells.Item(7, 28).Value = 14.5
$ws.Cells.Item(7, 30).Value = 18
$ws.Cells.Item(8, 6).Value = 2.36
$ws.Cells.Item(8, 7).Value = 2.6
$ws.Cells.Item(8, 9).Value = 3.1
$ws.Cells.Item(8, 12).Value = 1.18
$ws.Cells.Item(8, 14).Value = 7.2
$ws.Cells.Item(8, 16).Value = 3.05
$ws.Cells.Item(8, 17).Value = 1.38
$ws.Cells.Item(8, 18).Value = 1.86
$ws.Cells.Item(8, 19).Value = 1.92
$ws.Cells.Item(8, 20).Value = 1.39
$ws.Cells.Item(8, 21).Value = 3
$ws.Cells.Item(8, 22).Value = 1.48
$ws.Cells.Item(8, 23).Value = 1.62
$ws.Cells.Item(8, 24).Value = 980
$ws.Cells.Item(8, 25).Value = 980
$ws.Cells.Item(8, 26).Value = 980
$ws.Cells.Item(8, 27).Value = 55
$ws.Cells.Item(8, 28).Value = 980
$ws.Cells.Item(8, 29).Value = 14
$ws.Cells.Item(8, 31).Value = 980
$ws.Cells.Item(8, 32).Value = 980
$ws.Cells.Item(8, 33).Value = 16
$ws.Cells.Item(8, 34).Value = 16.5
$ws.Cells.Item(8, 35).Value = 980
$ws.Cells.Item(8, 36).Value = 980
$ws.Cells.Item(8, 37).Value = 980
$ws.Cells.Item(8, 38).Value = 980
$ws.Cells.Item(8, 39).Value = 980
$ws.Cells.Item(8, 40).Value = 11.5
$ws.Cells.Item(8, 41).Value = 14.5
$ws.Cells.Item(9, 29).Value = 980
$ws.Cells.Item(10, 14).Value = 6.6
$ws.Cells.Item(10, 15).Value = 1.14
$ws.Cells.Item(10, 16).Value = 2.9
$ws.Cells.Item(10, 17).Value = 1.41
$ws.Cells.Item(10, 18).Value = 1.77
$ws.Cells.Item(10, 19).Value = 2.04
$ws.Cells.Item(10, 20).Value = 1.47
$ws.Cells.Item(10, 21).Value = 2.66
$ws.Cells.Item(10, 24).Value = 42
$ws.Cells.Item(10, 25).Value = 32
$ws.Cells.Item(10, 26).Value = 46
$ws.Cells.Item(10, 27).Value = 90
$ws.Cells.Item(10, 28).Value = 20
$ws.Cells.Item(10, 29).Value = 14.5
$ws.Cells.Item(10, 30).Value = 22
$ws.Cells.Item(10, 31).Value = 46
$ws.Cells.Item(10, 32).Value = 20
$ws.Cells.Item(10, 33).Value = 14
$ws.Cells.Item(10, 34).Value = 19
$ws.Cells.Item(10, 35).Value = 44
$ws.Cells.Item(10, 36).Value = 27
$ws.Cells.Item(10, 37).Value = 21
$ws.Cells.Item(10, 38).Value = 29
$ws.Cells.Item(10, 39).Value = 60
$ws.Cells.Item(10, 40).Value = 8.199999999999999
$ws.Cells.Item(10, 41).Value = 28
$ws.Cells.Item(11, 8).Value = 3.75
$ws.Cells.Item(11, 19).Value = 3.75
$ws.Cells.Item(12, 12).Value = 1.53
$ws.Cells.Item(12, 41).Value = 44
$ws.Cells.Item(13, 9).Value = 15.5
$ws.Cells.Item(13, 17).Value = 1.7
